$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Each (cell, new value) pair below mirrors one <c> text change in the
# authoritative OOXML diff. NumberFormat is forced to "@" (Text) before the
# assignment so Excel does not reinterpret dotted/percent strings as numbers,
# then the style is reset to "Normal" so the saved cell keeps the workbook's
# original (default) style index -- only the cell VALUE changes, as in the diff.

$updates = @(
    @("D2", "27.131.50"),
    @("E2", "  -1.42%  "),
    @("D3", "1.799.21"),
    @("E3", "  -2.20%  "),
    @("D4", "1.008"),
    @("E4", "  +0.49%  "),
    @("D5", "1.008"),
    @("D6", "308.33"),
    @("E6", "  -1.62%  "),
    @("D7", "0.4177"),
    @("E7", "  -1.70%  "),
    @("D8", "0.3555"),
    @("E8", "  -2.97%  "),
    @("D9", "0.07056"),
    @("E9", "  -2.79%  "),
    @("D10", "0.8433"),
    @("E10", "  -2.93%  "),
    @("D11", "2.010.04"),
    @("E11", "  +8.49%  "),
    @("D12", "20.19"),
    @("E12", "  -2.87%  "),
    @("D13", "5.277"),
    @("E13", "  -2.23%  "),
    @("D14", "6.338"),
    @("E14", "  -2.81%  "),
    @("D15", "0.06788"),
    @("E15", "  -2.22%  "),
    @("D16", "1.011"),
    @("E16", "  +0.64%  "),
    @("D17", "79.89"),
    @("E17", "  -0.43%  "),
    @("D18", "0.000008710"),
    @("E18", "  -3.42%  "),
    @("E19", "  +0.58%  "),
    @("D20", "15.03"),
    @("E20", "  -2.74%  "),
    @("D21", "27.264.31"),
    @("E21", "  -1.17%  "),
    @("D22", "5.042"),
    @("E22", "  -0.28%  "),
    @("E23", "  -0.84%  "),
    @("D24", "2.105.52"),
    @("E24", "  +1.53%  "),
    @("D25", "1.949"),
    @("E25", "  -0.39%  "),
    @("D26", "153.16"),
    @("D27", "18.11"),
    @("E27", "  -1.19%  "),
    @("D28", "5.028"),
    @("E28", "  -4.10%  "),
    @("D29", "112.54"),
    @("E29", "  -2.38%  "),
    @("D30", "1.654"),
    @("E30", "  -10.83%  "),
    @("D31", "0.08872"),
    @("E31", "  +0.01%  "),
    @("B32", "HuobiToken"),
    @("C32", "https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht"),
    @("D32", "2.876"),
    @("E32", "  -2.73%  "),
    @("B33", "ImmutableX"),
    @("C33", "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"),
    @("D33", "0.7183"),
    @("E33", "  -7.31%  "),
    @("D34", "4.341"),
    @("E34", "  -4.46%  "),
    @("E35", "  +0.55%  "),
    @("D36", "1.077"),
    @("E36", "  -6.53%  "),
    @("D37", "1.076"),
    @("E37", "  -2.28%  "),
    @("D38", "0.01892"),
    @("E38", "  -2.78%  "),
    @("D39", "0.05092"),
    @("E39", "  -5.52%  "),
    @("E40", "  -2.91%  "),
    @("E41", "  -3.65%  "),
    @("D42", "2.608"),
    @("E42", "  -7.60%  "),
    @("D43", "6.168"),
    @("E43", "  -8.49%  "),
    @("D44", "8.067"),
    @("E44", "  -5.34%  "),
    @("D45", "1.008"),
    @("E45", "  +0.62%  "),
    @("D46", "104.37"),
    @("E46", "  -1.96%  "),
    @("D47", "10.17"),
    @("E47", "  -3.60%  "),
    @("D48", "0.06314"),
    @("E48", "  -3.27%  "),
    @("E49", "  -4.34%  "),
    @("E50", "  -3.12%  "),
    @("D51", "62.58"),
    @("E51", "  -2.79%  "),
)

foreach ($u in $updates) {
    $cell = $ws.Range($u[0])
    $cell.NumberFormat = "@"
    $cell.Value = $u[1]
    $cell.Style = "Normal"
}
